$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("W2019PHL")

$newRows = @(
    @("DOX_ND30", "ABX_DISK"),
    @("DOX_NM",   "ABX_MIC"),
    @("SSS_ND200","ABX_DISK"),
    @("SSS_NM",   "ABX_MIC")
)

$startRow = 169
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$ws.Range("B172").Select()
